$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values, forcing text storage to preserve exact formatting ---
$priceUpdates = @{
    "D2" = "40.034.46"
    "D3" = "2.344.51"
    "D5" = "310.29"
    "D6" = "84.72"
    "D7" = "0.525"
    "D9" = "0.481"
    "D10" = "0.0805"
    "D11" = "30.00"
    "D13" = "2.707.75"
    "D14" = "6.40"
    "D15" = "14.82"
    "D16" = "2.357.40"
    "D17" = "0.758"
    "D18" = "40.045.81"
    "D20" = "6.08"
    "D21" = "68.12"
    "D22" = "10.58"
    "D23" = "234.80"
    "D24" = "2.54"
    "D27" = "23.59"
    "D29" = "9.23"
    "D30" = "34.93"
    "D31" = "153.61"
    "D33" = "5.08"
    "D35" = "0.0716"
    "D36" = "0.113"
    "D37" = "2.77"
    "D38" = "0.0989"
    "D39" = "15.53"
    "D42" = "1.967.62"
    "D43" = "2.25"
    "D45" = "17.41"
    "D46" = "9.43"
    "D47" = "2.67"
    "D48" = "2.568.27"
    "D49" = "92.97"
    "D50" = "70.18"
    "D51" = "49.97"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# --- Update Volume(1h) (column E) values ---
$volumeUpdates = @{
    "E2" = "  -2.97%  "
    "E3" = "  -3.76%  "
    "E4" = "  -0.02%  "
    "E5" = "  -2.53%  "
    "E6" = "  -5.71%  "
    "E7" = "  -2.75%  "
    "E8" = "  -0.04%  "
    "E9" = "  -3.92%  "
    "E10" = "  -4.15%  "
    "E11" = "  -6.74%  "
    "E12" = "  +0.51%  "
    "E13" = "  -3.60%  "
    "E14" = "  -5.15%  "
    "E15" = "  -5.61%  "
    "E16" = "  -2.35%  "
    "E17" = "  -2.61%  "
    "E18" = "  -2.71%  "
    "E19" = "  -3.14%  "
    "E20" = "  -3.51%  "
    "E21" = "  -5.37%  "
    "E22" = "  -4.88%  "
    "E23" = "  -0.79%  "
    "E24" = "  -5.97%  "
    "E25" = "  -0.05%  "
    "E27" = "  -2.77%  "
    "E28" = "  -4.47%  "
    "E29" = "  -4.20%  "
    "E30" = "  +0.37%  "
    "E31" = "  -1.71%  "
    "E32" = "  -0.01%  "
    "E33" = "  -3.99%  "
    "E34" = "  -0.96%  "
    "E35" = "  -4.39%  "
    "E36" = "  -1.11%  "
    "E37" = "  -6.59%  "
    "E38" = "  -2.61%  "
    "E39" = "  -7.90%  "
    "E40" = "  -4.12%  "
    "E41" = "  -1.87%  "
    "E42" = "  -1.80%  "
    "E43" = "  +1.35%  "
    "E44" = "  -4.37%  "
    "E45" = "  -7.47%  "
    "E46" = "  -1.38%  "
    "E47" = "  -8.78%  "
    "E48" = "  -3.66%  "
    "E49" = "  -2.33%  "
    "E50" = "  -5.01%  "
    "E51" = "  -4.16%  "
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
